$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(10, 7).Value = 3
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(12, 7).Value = 2
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(17, 7).Value = 1
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(19, 7).Value = 1
$ws.Cells.Item(20, 7).Value = 1
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(27, 7).Value = 2
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(31, 7).Value = 3
$ws.Cells.Item(32, 7).Value = 1
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(34, 7).Value = 2
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(37, 7).Value = 3
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(41, 7).Value = 2
$ws.Cells.Item(42, 7).Value = 2
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(45, 7).Value = 1
$ws.Cells.Item(46, 7).Value = 2
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(48, 7).Value = 1
$ws.Cells.Item(49, 7).Value = 2
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(53, 7).Value = 2
$ws.Cells.Item(54, 7).Value = 2
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(57, 7).Value = 1
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(59, 7).Value = 2
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(65, 7).Value = 2
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(68, 7).Value = 2
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(72, 7).Value = 2
